$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Add a new row to the table (this extends the table range, the
# worksheet dimension and the autofilter range from A1:J92 to A1:J93).
$newRow = $tbl.ListRows.Add()

# Copy the formatting (number format, borders, fill, font, alignment)
# of the previous data row onto the freshly added row so the new row
# keeps looking like the rest of the table.
$fmtSrc = $ws.Range("A89:J89")
$fmtDst = $ws.Range("A93:J93")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122) # xlPasteFormats

# Fill in the new row's values.
$r = $newRow.Range
$r.Item(1).Value = 43993
$r.Item(2).Value = 86328
$r.Item(3).Value = 702
$r.Item(4).Value = 1490
$r.Item(5).Value = 2
$r.Item(6).Value = 6
$r.Item(7).Value = 0
$r.Item(8).Value = 0
$r.Item(9).Value = 109
$r.Item(10).Value = 0

# Match the author's selection state after entering the new row.
[void]$ws.Range("A93:J93").Select()
